# "Product Definition feedback with block diagram start"
#
# Hardware Development Process sheet: the "Estimate Architecture Task
# Hours" / "Estimate Task Hours" placeholder rows (12 & 13) are turned
# into two real deliverable rows - "Eagle Setup" and "Read through Eagle
# Tutorials" - each with hours/start/due/time data, and a new blank
# spacer row is inserted right after them (pushing everything else down
# by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Hardware Development Process" (tabSelected sheet)

# --- Row 13: "Eagle Setup" -------------------------------------------
$ws.Range("A13").Value = "Eagle Setup"

# --- Row 12: "Read through Eagle Tutorials " --------------------------
$ws.Range("A12").Value = "Read through Eagle Tutorials "

# Hours column (B) - copy the numeric/centered format used by the rows
# above and fill in the "hours" values.
$ws.Range("B6").Copy()
$ws.Range("B12:B13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B12").Value = 3
$ws.Range("B13").Value = 3

# Start/Due date columns (C/D) - copy the date format used by the rows
# above and fill in the dates.
$ws.Range("C6:D6").Copy()
$ws.Range("C12:D13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C12").Value = 42885
$ws.Range("D12").Value = 42891
$ws.Range("C13").Value = 42885
$ws.Range("D13").Value = 42891

# Time column (E)
$ws.Range("E12").Value = "2:00pm"
$ws.Range("E13").Value = "2:00pm"

# --- New blank spacer row, pushing the rest of the table down --------
$ws.Rows(14).Insert()

# --- Selection left on A23 after the edit -----------------------------
$ws.Range("A23").Select() | Out-Null
